$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G2").Value = "2016-03-03 11:05:25"
$zhcn.Range("G3").Value = "2016-03-03 11:05:25"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G2").Value = "2016-03-03 11:05:47"
$dede.Range("G3").Value = "2016-03-03 11:05:47"
